$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns that vary row-to-row (D, L, M, N, O, P, R, S) are being
# reshuffled: each data row (2-30) ends up holding the values that
# previously belonged to a different row. Capture the "movable" columns
# for every source row first, then write them back out in the new
# arrangement so we never read a cell after it has already been
# overwritten.

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value()
    }
    $snapshot[$r] = $rowData
}

# New row number -> row number whose values it should receive.
$mapping = @(
    @{New=2;  Old=29},
    @{New=3;  Old=18},
    @{New=4;  Old=10},
    @{New=5;  Old=11},
    @{New=6;  Old=8},
    @{New=7;  Old=28},
    @{New=8;  Old=25},
    @{New=9;  Old=26},
    @{New=10; Old=27},
    @{New=11; Old=21},
    @{New=12; Old=22},
    @{New=13; Old=19},
    @{New=14; Old=20},
    @{New=15; Old=2},
    @{New=16; Old=3},
    @{New=17; Old=4},
    @{New=18; Old=30},
    @{New=19; Old=12},
    @{New=20; Old=13},
    @{New=21; Old=5},
    @{New=22; Old=6},
    @{New=23; Old=7},
    @{New=24; Old=23},
    @{New=25; Old=14},
    @{New=26; Old=15},
    @{New=27; Old=16},
    @{New=28; Old=17},
    @{New=29; Old=9},
    @{New=30; Old=24}
)

foreach ($entry in $mapping) {
    $newRow = $entry.New
    $oldRow = $entry.Old
    $rowData = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $rowData[$col]
    }
}
